$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Update the time_taken values (column F) in the "data" sheet
$dataSheet.Range("F2").Value = "2021-10-05 14:33:39.342872"
$dataSheet.Range("F3").Value = "2021-10-05 14:33:39.342880"
$dataSheet.Range("F4").Value = "2021-10-05 14:33:39.342883"
$dataSheet.Range("F5").Value = "2021-10-05 14:33:39.342885"
$dataSheet.Range("F6").Value = "2021-10-05 14:33:39.342888"
$dataSheet.Range("F7").Value = "2021-10-05 14:33:39.342891"
$dataSheet.Range("F8").Value = "2021-10-05 14:33:39.342893"
$dataSheet.Range("F9").Value = "2021-10-05 14:33:39.342896"
$dataSheet.Range("F10").Value = "2021-10-05 14:33:39.342899"
$dataSheet.Range("F11").Value = "2021-10-05 14:33:39.342901"
$dataSheet.Range("F12").Value = "2021-10-05 14:33:39.342904"
$dataSheet.Range("F13").Value = "2021-10-05 14:33:39.342906"
$dataSheet.Range("F14").Value = "2021-10-05 14:33:39.342909"
$dataSheet.Range("F15").Value = "2021-10-05 14:33:39.342911"
$dataSheet.Range("F16").Value = "2021-10-05 14:33:39.342914"

# Add a new "metadata" worksheet right after the "data" sheet
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$ws2.Name = "metadata"

# Copy cell formatting (bold header style, bordered A2 style) from "data" sheet
$dataSheet.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)

$dataSheet.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

# Header row values
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row values
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Cutis Laxa"
$ws2.Range("C2").Value = 3129
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "0.8"
$ws2.Range("E2").Value = "2021-07-22T23:49:12.737975Z"
$ws2.Range("F2").Value = "2021-10-05 14:33:39.338975"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3129/?format=json"
